# 2024-12-11, Dewey artikel update
# Adds a new bibliography-style paragraph ("Edel, A. and Flower, E. , Introduction,
# LW7, p. vii- xxxv") right before the last (empty) paragraph of the document,
# matching the formatting (Courier New, 360 auto line spacing, en-US) already used
# by the rest of the reference list, including the proofErr spell/grammar markers
# Word leaves around "Edel" and "E. ,".

$d = $word.ActiveDocument

# The document ends with a lone empty paragraph (just a paragraph mark) right
# before the section properties. Insert a brand-new paragraph right before it,
# which keeps that trailing empty paragraph exactly where/what it was.
$trailing = $d.Paragraphs($d.Paragraphs.Count)
$trailing.Range.InsertParagraphBefore() | Out-Null

# The freshly inserted (still empty) paragraph is now second-to-last.
$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$rangeStart = $newPara.Range.Start
$rangeEnd = $newPara.Range.End

# Replace the new paragraph's content (its paragraph mark included) with the
# full Open XML for the paragraph: its pPr plus the four runs and the
# spellcheck/grammar-check proofErr bookmarks around "Edel" and "E. ,", exactly
# as Word itself recorded them.
$target = $d.Range($rangeStart, $rangeEnd)

$xmlFrag = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:spacing w:line="360" w:lineRule="auto"/>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
            </w:pPr>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>Edel</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve">, A. and Flower, </w:t>
            </w:r>
            <w:proofErr w:type="gramStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t>E. ,</w:t>
            </w:r>
            <w:proofErr w:type="gramEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
                <w:lang w:val="en-US"/>
              </w:rPr>
              <w:t xml:space="preserve"> Introduction, LW7, p. vii- xxxv</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.InsertXML($xmlFrag) | Out-Null
